# edit.ps1 - applies the USERSTORIES.docx revision described by the diff.
$d = $word.ActiveDocument

function Replace-Text($range, $old, $new) {
    $ok = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace failed for: $old"
    }
}

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the document (just before the
#    "empty sz20/222222" paragraph near the Database Manager section) to the
#    empty centered paragraph right under the "USER STORIES" heading.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$d.Bookmarks.Add("_GoBack", $p2.Range) | Out-Null

# ---------------------------------------------------------------------------
# 2. Player user stories 1-6 : add a comma after "As a Player" and update the
#    remaining wording.
# ---------------------------------------------------------------------------
Replace-Text $d.Paragraphs.Item(4).Range `
    "1. As a Player I want to be able to view my HP/Mana/Food so that I know what my options are." `
    "1. As a Player, I want to be able to view my HP/Mana/Food so that I can decide what actions to take."

Replace-Text $d.Paragraphs.Item(5).Range `
    "2. As a Player I want to be able to log in so that my score can be saved across games." `
    "2. As a Player, I want to be able to log in with a unique username so that my score can be saved across games and so that my data is not mixed up with that of other players."

Replace-Text $d.Paragraphs.Item(6).Range `
    "3. As a Player I want to be able to EAT so that I can heal." `
    "3. As a Player, I want to be able to EAT so that I can heal."

Replace-Text $d.Paragraphs.Item(7).Range `
    "4. As a Player I want to be able to ATTACK so that I can damage the enemy." `
    "4. As a Player, I want to be able to ATTACK so that I can damage the enemy."

Replace-Text $d.Paragraphs.Item(8).Range `
    "5. As a Player I want to be able to use SP ATTACK so that I can damage the enemy." `
    "5. As a Player, I want to be able to use SP ATTACK so that I can damage the enemy."

Replace-Text $d.Paragraphs.Item(9).Range `
    "6. As a Player I want to be able to view my enemies HP/Mana/Food so that I know when to attack." `
    ("6. As a Player, I want to be able to view my opponent" + [char]0x2019 + "s HP/Mana/Food so that I can decide what actions to take.")

# ---------------------------------------------------------------------------
# 3. Story 7 : add a comma after "As a Player".
# ---------------------------------------------------------------------------
Replace-Text $d.Paragraphs.Item(10).Range `
    "7. As a Player I want to be able to see my score so that I know what my score is." `
    "7. As a Player, I want to be able to see my score so that I know what my score is."

# ---------------------------------------------------------------------------
# 4. Story 8 : add a comma after "As a Player" (keep the rest / formatting).
# ---------------------------------------------------------------------------
Replace-Text $d.Paragraphs.Item(11).Range `
    "8. As a Player I want to be able to navigate" `
    "8. As a Player, I want to be able to navigate"

# ---------------------------------------------------------------------------
# 5. Insert new story 11 ("view the result of dice rolls...") before the
#    existing "11. As a player, I want to have different options..." item,
#    which will become story 12 (and every later item shifts by one).
# ---------------------------------------------------------------------------
$storyDifferentOptions = $d.Paragraphs.Item(14)
$storyDifferentOptions.Range.InsertParagraphBefore()
$newStory11 = $d.Paragraphs.Item(14)
$newStory11.Range.Text = "11. As a player, I want to be able to view the result of dice rolls so that I know how the results of actions are computed."

# Renumber the following items (12-15 -> 12 stays same value but need to bump
# since old 11 -> 12, old 12 -> 13, old 13 -> 14, old 14 -> 15).
Replace-Text $d.Paragraphs.Item(15).Range "11." "12."
Replace-Text $d.Paragraphs.Item(16).Range "12." "13."
Replace-Text $d.Paragraphs.Item(17).Range "13." "14."

# Story (old 14, now 15) : renumber and change "view past" -> "view saved".
Replace-Text $d.Paragraphs.Item(18).Range "14." "15."
Replace-Text $d.Paragraphs.Item(18).Range "view past statistics" "view saved statistics"

# ---------------------------------------------------------------------------
# 6. Insert a new blank paragraph (matching the spacing/formatting of the
#    paragraph mark of story 15) right after it.
# ---------------------------------------------------------------------------
$story15 = $d.Paragraphs.Item(18)
$story15.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 7. Old story 15 ("exit during a battle...") becomes story 16; remove the
#    sz/szCs from the paragraph mark and renumber "15." -> "16.".
# ---------------------------------------------------------------------------
$story16 = $d.Paragraphs.Item(20)
Replace-Text $story16.Range "15." "16."
$story16.Range.Font.Size = 11

# ---------------------------------------------------------------------------
# 8. Add a brand-new story 17 ("save statistics from battles...") into what
#    used to be an empty paragraph right after story 16, then append three
#    blank paragraphs (replacing the old tab-stop paragraphs) before the
#    "Database Manager" heading.
# ---------------------------------------------------------------------------
$emptyAfter16 = $d.Paragraphs.Item(21)
$emptyAfter16.Range.Text = "17. As a player, I want to be able to save statistics from battles so that I can view them when I want to."

$dbManagerPara = $d.Paragraphs.Item(23)
$dbManagerPara.Range.InsertParagraphBefore()
$dbManagerPara.Range.InsertParagraphBefore()

# Remove the tab characters / tab stops from what used to be the two
# tab-stop paragraphs, leaving plain empty paragraphs.
$tabPara1 = $d.Paragraphs.Item(24)
$tabPara1.Range.Text = ""
$tabPara1.Format.TabStops.ClearAll()

$tabPara2 = $d.Paragraphs.Item(25)
$tabPara2.Format.TabStops.ClearAll()

# ---------------------------------------------------------------------------
# 9. "Database Manager" heading gets a lastRenderedPageBreak marker before
#    its text run.
# ---------------------------------------------------------------------------
$dbHeading = $d.Paragraphs.Item(26)
$dbHeading.Range.InsertBefore([char]2)
